$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 32: Antonio ADM ---
$ws.Range("A32").Value = "Antonio ADM"
$ws.Range("A32").Borders.LineStyle = 1
$ws.Range("A32").Borders.Color = 0

$ws.Range("B32").Value = 554384356465
$ws.Range("B32").NumberFormat = "0"
$ws.Range("B32").HorizontalAlignment = -4108
$ws.Range("B32").Borders.LineStyle = 1
$ws.Range("B32").Borders.Color = 0

$ws.Range("C32").Value = 554384356465
$ws.Range("C32").NumberFormat = "0"
$ws.Range("C32").HorizontalAlignment = -4108
$ws.Range("C32").Borders.LineStyle = 1
$ws.Range("C32").Borders.Color = 0

# --- Row 33: Pedro ADM ---
$ws.Range("A33").Value = "Pedro ADM"
$ws.Range("A33").Borders.LineStyle = 1
$ws.Range("A33").Borders.Color = 0

$ws.Range("B33").Value = 5543996440402
$ws.Range("B33").NumberFormat = "0"
$ws.Range("B33").HorizontalAlignment = -4108
$ws.Range("B33").Borders.LineStyle = 1
$ws.Range("B33").Borders.Color = 0

$ws.Range("C33").Value = 5543996440402
$ws.Range("C33").NumberFormat = "0"
$ws.Range("C33").HorizontalAlignment = -4108
$ws.Range("C33").Borders.LineStyle = 1
$ws.Range("C33").Borders.Color = 0
